$wb = $excel.ActiveWorkbook

# New column C values (Reaction_number) for rows 2..16 of each sheet, in order.
# Column B (Cutoff index) for data row i (0-based, i = 0..14) becomes i + 5.
$newValues = @{
    "NBR" = @(740, 721, 716, 715, 724, 715, 704, 704, 701, 692, 686, 680, 681, 651, 0)
    "BAR" = @(752, 755, 767, 750, 746, 747, 751, 747, 752, 750, 744, 754, 752, 755, 0)
}

foreach ($sheetName in @("NBR", "BAR")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $newValues[$sheetName]

    # Update rows 2..16 (15 data rows) with shifted B values and new C values.
    for ($i = 0; $i -lt 15; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $i + 5
        $ws.Cells.Item($row, 3).Value = $values[$i]
    }

    # Remove the now-extra rows 17..20 (the table shrank from 19 to 15 entries).
    $ws.Rows("17:20").Delete()
}
